$d = $word.ActiveDocument

# Helper: within the text matched by $searchPhrase (found via Find.Execute),
# replace the single character at zero-based $offsetInPhrase with $newChar.
#
# A direct Range.Text assignment on that character would make the engine
# recombine it with its (identically formatted) neighbouring runs into one
# run. To match the target document - where the original run is split into
# two sibling runs with identical rPr ("6"+"2", "4"+"4") - we briefly toggle
# a character property (Bold) on just that character before editing it, and
# clear it again afterwards. That keeps the edited character as its own run
# instead of being re-merged with the runs on either side.
function Split-Digit($searchPhrase, $offsetInPhrase, $newChar) {
    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($searchPhrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $searchPhrase"
        return
    }

    $target = $d.Range($rng.Start + $offsetInPhrase, $rng.Start + $offsetInPhrase + 1)
    $target.Font.Bold = $true
    $target.Text = $newChar

    $newRun = $d.Range($rng.Start + $offsetInPhrase, $rng.Start + $offsetInPhrase + 1)
    $newRun.Font.Bold = $false
}

# "Florent : 60h" -> "Florent : 62h" (split the "60" run into "6" + "2")
Split-Digit "Florent : 60h" 11 "2"

# "Alexis : 41h" -> "Alexis : 44h" (split the "41" run into "4" + "4")
Split-Digit "Alexis : 41h" 10 "4"
